$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.147582333333333
$ws.Range("H2").Value = 3.442747
$ws.Range("I2").Value = 0.0006681225322425688
$ws.Range("J2").Value = 0.0006681225322425686
$ws.Range("M2").Value = 1.186139
$ws.Range("N2").Value = 3.558416999999999
$ws.Range("O2").Value = 0.07586413741082078
$ws.Range("P2").Value = 0.07586413741082078
$ws.Range("Q2").Value = 1.361192161277666
$ws.Range("R2").Value = 12.250729451499
$ws.Range("S2").Value = [double]"5.068653959331577E-05"
$ws.Range("T2").Value = [double]"5.068653959331576E-05"
$ws.Range("G3").Value = 1.147582333333333
$ws.Range("H3").Value = 3.442747
$ws.Range("I3").Value = 0.0006681225322425688
$ws.Range("J3").Value = 0.0006681225322425686
$ws.Range("O3").Value = 0.7473293058134135
$ws.Range("P3").Value = 0.7473293058134136
$ws.Range("Q3").Value = 13.40895484591911
$ws.Range("R3").Value = 120.680593613272
$ws.Range("S3").Value = 0.0004993075482191389
$ws.Range("T3").Value = 0.0004993075482191389
$ws.Range("G4").Value = 1.147582333333333
$ws.Range("H4").Value = 3.442747
$ws.Range("I4").Value = 0.0006681225322425688
$ws.Range("J4").Value = 0.0006681225322425686
$ws.Range("M4").Value = 2.764377999999999
$ws.Range("N4").Value = 8.293133999999998
$ws.Range("O4").Value = 0.1768065567757656
$ws.Range("P4").Value = 0.1768065567757657
$ws.Range("Q4").Value = 3.172351355455332
$ws.Range("R4").Value = 28.55116219909799
$ws.Range("S4").Value = 0.0001181284444301141
$ws.Range("T4").Value = 0.0001181284444301141
$ws.Range("I5").Value = 0.9769430428898132
$ws.Range("J5").Value = 0.9769430428898132
$ws.Range("M5").Value = 1.186139
$ws.Range("N5").Value = 3.558416999999999
$ws.Range("O5").Value = 0.07586413741082078
$ws.Range("P5").Value = 0.07586413741082078
$ws.Range("Q5").Value = 1990.36426377185
$ws.Range("R5").Value = 17913.27837394665
$ws.Range("S5").Value = 0.07411494124833816
$ws.Range("T5").Value = 0.07411494124833816
$ws.Range("I6").Value = 0.9769430428898132
$ws.Range("J6").Value = 0.9769430428898132
$ws.Range("O6").Value = 0.7473293058134135
$ws.Range("P6").Value = 0.7473293058134136
$ws.Range("S6").Value = 0.7300981660620879
$ws.Range("T6").Value = 0.730098166062088
$ws.Range("I7").Value = 0.9769430428898132
$ws.Range("J7").Value = 0.9769430428898132
$ws.Range("M7").Value = 2.764377999999999
$ws.Range("N7").Value = 8.293133999999998
$ws.Range("O7").Value = 0.1768065567757656
$ws.Range("P7").Value = 0.1768065567757657
$ws.Range("Q7").Value = 4638.679937812601
$ws.Range("S7").Value = 0.172729935579387
$ws.Range("T7").Value = 0.172729935579387
$ws.Range("I8").Value = 0.02238883457794425
$ws.Range("J8").Value = 0.02238883457794424
$ws.Range("M8").Value = 1.186139
$ws.Range("N8").Value = 3.558416999999999
$ws.Range("O8").Value = 0.07586413741082078
$ws.Range("P8").Value = 0.07586413741082078
$ws.Range("Q8").Value = 45.61364818119266
$ws.Range("R8").Value = 410.5228336307339
$ws.Range("S8").Value = 0.001698509622889298
$ws.Range("T8").Value = 0.001698509622889298
$ws.Range("I9").Value = 0.02238883457794425
$ws.Range("J9").Value = 0.02238883457794424
$ws.Range("O9").Value = 0.7473293058134135
$ws.Range("P9").Value = 0.7473293058134136
$ws.Range("Q9").Value = 449.3350507140391
$ws.Range("S9").Value = 0.01673183220310642
$ws.Range("T9").Value = 0.01673183220310642
$ws.Range("I10").Value = 0.02238883457794425
$ws.Range("J10").Value = 0.02238883457794424
$ws.Range("M10").Value = 2.764377999999999
$ws.Range("N10").Value = 8.293133999999998
$ws.Range("O10").Value = 0.1768065567757656
$ws.Range("P10").Value = 0.1768065567757657
$ws.Range("R10").Value = 956.7515188240678
$ws.Range("S10").Value = 0.003958492751948525
$ws.Range("T10").Value = 0.003958492751948525
